$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1. Insert the two new sheets.
#    Target sheetId order needs: SignUp(3) "ManageListings(5) SignIn(1) ShareSkill(4) Profile(2)
#    sheetId is assigned sequentially as sheets are *created*, so ShareSkill
#    (which ends up with the lower new id, 4) must be created before
#    "ManageListings (which ends up with id 5), even though ManageListings is
#    positioned earlier in the tab strip. We fix up the tab order afterwards.
# ---------------------------------------------------------------------------

# Create ShareSkill immediately before Profile (gets sheetId 4).
$profile = $wb.Worksheets.Item("Profile")
$shareSkill = $wb.Worksheets.Add($null, $profile)
$shareSkill.Name = "ShareSkill"

# Create "ManageListings immediately after SignUp (gets sheetId 5).
$signUp = $wb.Worksheets.Item("SignUp")
$manageListings = $wb.Worksheets.Add($null, $signUp)
$manageListings.Name = "`"ManageListings"

# Re-fetch sheets by name (index-based refs go stale after inserts) and move
# ShareSkill so it sits right before Profile, after SignIn, giving:
#   SignUp, "ManageListings, SignIn, ShareSkill, Profile
$shareSkill = $wb.Worksheets.Item("ShareSkill")
$profile = $wb.Worksheets.Item("Profile")
$shareSkill.Move($profile)

# ---------------------------------------------------------------------------
# 2. SignIn's selection memory moves from C2 to C3 (no longer the active tab;
#    "ManageListings becomes the active/selected tab instead).
# ---------------------------------------------------------------------------
$signIn = $wb.Worksheets.Item("SignIn")
$signIn.Activate()
$signIn.Range("C3").Select()

# ---------------------------------------------------------------------------
# 3. Seed "ManageListings header row first so the "Title" shared string is
#    introduced before ShareSkill needs it too (matches shared-string order).
# ---------------------------------------------------------------------------
$manageListings = $wb.Worksheets.Item("`"ManageListings")
$manageListings.Range("A1").Value = "Title"

# ---------------------------------------------------------------------------
# 4. Populate ShareSkill fully (row 1 headers, row 2 data).
# ---------------------------------------------------------------------------
$shareSkill = $wb.Worksheets.Item("ShareSkill")

$shareSkill.Range("A1").Value = "Title"
$shareSkill.Range("B1").Value = "Description"
$shareSkill.Range("C1").Value = "Category"
$shareSkill.Range("D1").Value = "SubCategory"
$shareSkill.Range("E1").Value = "Tags"
$shareSkill.Range("F1").Value = "ServiceType"
$shareSkill.Range("G1").Value = "LocationType"
$shareSkill.Range("H1").Value = "Startdate"
$shareSkill.Range("I1").Value = "Enddate"
$shareSkill.Range("J1").Value = "Selectday"
$shareSkill.Range("K1").Value = "Starttime"
$shareSkill.Range("L1").Value = "Endtime"
$shareSkill.Range("M1").Value = "SkillTrade"
$shareSkill.Range("N1").Value = "Skill-Exchange"
$shareSkill.Range("O1").Value = "Credit"
$shareSkill.Range("P1").Value = "Active"

$shareSkill.Range("A2").Value = "Selenium"
$shareSkill.Range("B2").Value = "Would like to provide selenium training for beginners"
$shareSkill.Range("C2").Value = "Programming & Tech"
$shareSkill.Range("D2").Value = "QA"
$shareSkill.Range("E2").Value = "Testing"
$shareSkill.Range("F2").Value = "One-off service"
$shareSkill.Range("G2").Value = "On-site"

# Startdate / Enddate: numFmtId 14 ("mm-dd-yy"). Set format on H2 first, then
# clone that exact style onto I2 via a format-only paste so both cells share
# a single cellXf (setting NumberFormat independently on each cell would
# otherwise mint two separate, functionally-identical style entries).
$shareSkill.Range("H2").Value = 44663
$shareSkill.Range("H2").NumberFormat = "mm-dd-yy"
$shareSkill.Range("I2").Value = 44724
$shareSkill.Range("H2").Copy()
$shareSkill.Range("I2").PasteSpecial($xlPasteFormats)

$shareSkill.Range("J2").Value = "Mon"

# Starttime / Endtime: numFmtId 21 ("h:mm:ss"), same single-style trick.
$shareSkill.Range("K2").Value = 0.75
$shareSkill.Range("K2").NumberFormat = "h:mm:ss"
$shareSkill.Range("L2").Value = 0.83333333333333337
$shareSkill.Range("K2").Copy()
$shareSkill.Range("L2").PasteSpecial($xlPasteFormats)

$shareSkill.Range("M2").Value = "Skill-Exchange"
$shareSkill.Range("N2").Value = "Performance Testing"
$shareSkill.Range("O2").Value = "Credit"
$shareSkill.Range("P2").Value = "Hidden"

$wb.Application.CutCopyMode = $false

# Column widths / selection for ShareSkill.
$shareSkill.Columns.Item(2).ColumnWidth = 38.5
$shareSkill.Columns.Item(8).ColumnWidth = 15.333333333333334
$shareSkill.Columns.Item(9).ColumnWidth = 16.166666666666668
$shareSkill.Range("I2").Select()

# ---------------------------------------------------------------------------
# 5. Finish "ManageListings: header fill style (reuse the same style already
#    used by every other sheet's header row, via a format-only paste so no
#    new cellXf is minted), then the remaining data/strings.
# ---------------------------------------------------------------------------
$manageListings = $wb.Worksheets.Item("`"ManageListings")
$headerStyleSrc = $wb.Worksheets.Item("SignUp").Range("A1:B1")
$headerStyleSrc.Copy()
$manageListings.Range("A1:B1").PasteSpecial($xlPasteFormats)
$wb.Application.CutCopyMode = $false

$manageListings.Range("A2").Value = "Selenium"
$manageListings.Range("B1").Value = "Deleteaction"
$manageListings.Range("B2").Value = "Yes"

# Column widths / selection for "ManageListings.
$manageListings.Columns.Item(1).ColumnWidth = 12.5
$manageListings.Columns.Item(2).ColumnWidth = 23.833333333333332
$manageListings.Range("A1:B2").Select()
